# Apply the "top 50" re-pull with time filtering + two new columns (TVN, CTC).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two added columns
$ws.Range("F1").Value = "TVN"
$ws.Range("G1").Value = "CTC"

# Row 1
$ws.Range("A1").Value = "Index"
$ws.Range("B1").Value = "Time"
$ws.Range("C1").Value = "SegStart"
$ws.Range("D1").Value = "SegEnd"
$ws.Range("E1").Value = "AWC"

# Row 2
$ws.Range("A2").Value = 39
$ws.Range("B2").Value = "5:45 AM"
$ws.Range("C2").Value = 1169.48
$ws.Range("D2").Value = 1199.48
$ws.Range("E2").Value = 21.87
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2

# Row 3
$ws.Range("A3").Value = 50
$ws.Range("B3").Value = "5:50 AM"
$ws.Range("C3").Value = 1489.303333
$ws.Range("D3").Value = 1519.303333
$ws.Range("E3").Value = 2.51
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

# Row 4
$ws.Range("A4").Value = 84
$ws.Range("B4").Value = "6:07 AM"
$ws.Range("C4").Value = 2509.07
$ws.Range("D4").Value = 2539.07
$ws.Range("E4").Value = 3.404
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.8

# Row 5
$ws.Range("A5").Value = 93
$ws.Range("B5").Value = "6:11 AM"
$ws.Range("C5").Value = 2774.6
$ws.Range("D5").Value = 2804.6
$ws.Range("E5").Value = 16.77
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1

# Row 6
$ws.Range("A6").Value = 103
$ws.Range("B6").Value = "6:16 AM"
$ws.Range("C6").Value = 3079.71
$ws.Range("D6").Value = 3109.71
$ws.Range("E6").Value = 14.55
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1

# Row 7
$ws.Range("A7").Value = 116
$ws.Range("B7").Value = "6:23 AM"
$ws.Range("C7").Value = 3473.17
$ws.Range("D7").Value = 3503.17
$ws.Range("E7").Value = 1.063333
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0.333333

# Row 8
$ws.Range("A8").Value = 139
$ws.Range("B8").Value = "6:34 AM"
$ws.Range("C8").Value = 4149.392
$ws.Range("D8").Value = 4179.392
$ws.Range("E8").Value = 1.24
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0.2

# Row 9
$ws.Range("A9").Value = 146
$ws.Range("B9").Value = "6:38 AM"
$ws.Range("C9").Value = 4379.04
$ws.Range("D9").Value = 4409.04
$ws.Range("E9").Value = 3.49
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

# Row 10
$ws.Range("A10").Value = 170
$ws.Range("B10").Value = "6:50 AM"
$ws.Range("C10").Value = 5077.825
$ws.Range("D10").Value = 5107.825
$ws.Range("E10").Value = 0.4025
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0

# Row 11
$ws.Range("A11").Value = 179
$ws.Range("B11").Value = "6:54 AM"
$ws.Range("C11").Value = 5358.123333
$ws.Range("D11").Value = 5388.123333
$ws.Range("E11").Value = 1.47
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0

# Row 12
$ws.Range("A12").Value = 198
$ws.Range("B12").Value = "7:04 AM"
$ws.Range("C12").Value = 5933.15
$ws.Range("D12").Value = 5963.15
$ws.Range("E12").Value = 1.195
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0.5

# Row 13
$ws.Range("A13").Value = 214
$ws.Range("B13").Value = "7:12 AM"
$ws.Range("C13").Value = 6404.45
$ws.Range("D13").Value = 6434.45
$ws.Range("E13").Value = 1.525
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

# Row 14
$ws.Range("A14").Value = 221
$ws.Range("B14").Value = "7:16 AM"
$ws.Range("C14").Value = 6626.76
$ws.Range("D14").Value = 6656.76
$ws.Range("E14").Value = 3.56
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 1

# Row 15
$ws.Range("A15").Value = 250
$ws.Range("B15").Value = "7:30 AM"
$ws.Range("C15").Value = 7488.7825
$ws.Range("D15").Value = 7518.7825
$ws.Range("E15").Value = 0.08749999999999999
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0

# Row 16
$ws.Range("A16").Value = 277
$ws.Range("B16").Value = "7:43 AM"
$ws.Range("C16").Value = 8299.91
$ws.Range("D16").Value = 8329.91
$ws.Range("E16").Value = 4.62
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0

# Row 17
$ws.Range("A17").Value = 342
$ws.Range("B17").Value = "8:16 AM"
$ws.Range("C17").Value = 10238.38
$ws.Range("D17").Value = 10268.38
$ws.Range("E17").Value = 2.085
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0

# Row 18
$ws.Range("A18").Value = 387
$ws.Range("B18").Value = "8:39 AM"
$ws.Range("C18").Value = 11609.42
$ws.Range("D18").Value = 11639.42
$ws.Range("E18").Value = 8.06
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0

# Row 19
$ws.Range("A19").Value = 425
$ws.Range("B19").Value = "8:57 AM"
$ws.Range("C19").Value = 12725.265
$ws.Range("D19").Value = 12755.265
$ws.Range("E19").Value = 4.76
$ws.Range("F19").Value = 43.895
$ws.Range("G19").Value = 0

# Row 20
$ws.Range("A20").Value = 435
$ws.Range("B20").Value = "9:02 AM"
$ws.Range("C20").Value = 13043.06
$ws.Range("D20").Value = 13073.06
$ws.Range("E20").Value = 4.12
$ws.Range("F20").Value = 1.68
$ws.Range("G20").Value = 1

# Row 21
$ws.Range("A21").Value = 443
$ws.Range("B21").Value = "9:06 AM"
$ws.Range("C21").Value = 13276.226667
$ws.Range("D21").Value = 13306.226667
$ws.Range("E21").Value = 1.58
$ws.Range("F21").Value = 5.873333
$ws.Range("G21").Value = 0.333333

# Row 22
$ws.Range("A22").Value = 451
$ws.Range("B22").Value = "9:10 AM"
$ws.Range("C22").Value = 13511.86
$ws.Range("D22").Value = 13541.86
$ws.Range("E22").Value = 5.205
$ws.Range("F22").Value = 10.74
$ws.Range("G22").Value = 0.5

# Row 23
$ws.Range("A23").Value = 546
$ws.Range("B23").Value = "9:58 AM"
$ws.Range("C23").Value = 16367.32
$ws.Range("D23").Value = 16397.32
$ws.Range("E23").Value = 2.1925
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0

# Row 24
$ws.Range("A24").Value = 567
$ws.Range("B24").Value = "10:08 A"
$ws.Range("C24").Value = 16994.61
$ws.Range("D24").Value = 17024.61
$ws.Range("E24").Value = 7.28
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0

# Row 25
$ws.Range("A25").Value = 575
$ws.Range("B25").Value = "10:12 A"
$ws.Range("C25").Value = 17229.33
$ws.Range("D25").Value = 17259.33
$ws.Range("E25").Value = 1.03
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0

# Row 26
$ws.Range("A26").Value = 981
$ws.Range("B26").Value = "1:35 PM"
$ws.Range("C26").Value = 29415.15
$ws.Range("D26").Value = 29445.15
$ws.Range("E26").Value = 13.746667
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0

# Row 27
$ws.Range("A27").Value = 990
$ws.Range("B27").Value = "1:40 PM"
$ws.Range("C27").Value = 29686.38
$ws.Range("D27").Value = 29716.38
$ws.Range("E27").Value = 5.98
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0.5

# Row 28
$ws.Range("A28").Value = 1012
$ws.Range("B28").Value = "1:51 PM"
$ws.Range("C28").Value = 30341.53
$ws.Range("D28").Value = 30371.53
$ws.Range("E28").Value = 6.06
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0.25

# Row 29
$ws.Range("A29").Value = 1035
$ws.Range("B29").Value = "2:02 PM"
$ws.Range("C29").Value = 31029.245
$ws.Range("D29").Value = 31059.245
$ws.Range("E29").Value = 17.475
$ws.Range("F29").Value = 0.545
$ws.Range("G29").Value = 0

# Row 30
$ws.Range("A30").Value = 1048
$ws.Range("B30").Value = "2:09 PM"
$ws.Range("C30").Value = 31421.85
$ws.Range("D30").Value = 31451.85
$ws.Range("E30").Value = 4.695
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0

# Row 31
$ws.Range("A31").Value = 1104
$ws.Range("B31").Value = "2:37 PM"
$ws.Range("C31").Value = 33100.574
$ws.Range("D31").Value = 33130.574
$ws.Range("E31").Value = 2.998
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0

# Row 32
$ws.Range("A32").Value = 1141
$ws.Range("B32").Value = "2:55 PM"
$ws.Range("C32").Value = 34219.87
$ws.Range("D32").Value = 34249.87
$ws.Range("E32").Value = 21.48
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 2

# Row 33
$ws.Range("A33").Value = 1157
$ws.Range("B33").Value = "3:03 PM"
$ws.Range("C33").Value = 34686.29
$ws.Range("D33").Value = 34716.29
$ws.Range("E33").Value = 3
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0

# Row 34
$ws.Range("A34").Value = 1192
$ws.Range("B34").Value = "3:21 PM"
$ws.Range("C34").Value = 35742.3425
$ws.Range("D34").Value = 35772.3425
$ws.Range("E34").Value = 13.795
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0

# Row 35
$ws.Range("A35").Value = 1214
$ws.Range("B35").Value = "3:32 PM"
$ws.Range("C35").Value = 36391.61
$ws.Range("D35").Value = 36421.61
$ws.Range("E35").Value = 96.52
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 7

# Row 36
$ws.Range("A36").Value = 1224
$ws.Range("B36").Value = "3:37 PM"
$ws.Range("C36").Value = 36710.185
$ws.Range("D36").Value = 36740.185
$ws.Range("E36").Value = 3.49
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0.25

# Row 37
$ws.Range("A37").Value = 1242
$ws.Range("B37").Value = "3:46 PM"
$ws.Range("C37").Value = 37247.11
$ws.Range("D37").Value = 37277.11
$ws.Range("E37").Value = 5.1975
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0

# Row 38
$ws.Range("A38").Value = 1257
$ws.Range("B38").Value = "3:53 PM"
$ws.Range("C38").Value = 37684.8
$ws.Range("D38").Value = 37714.8
$ws.Range("E38").Value = 36.82
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 1

# Row 39
$ws.Range("A39").Value = 1267
$ws.Range("B39").Value = "3:58 PM"
$ws.Range("C39").Value = 37993.783333
$ws.Range("D39").Value = 38023.783333
$ws.Range("E39").Value = 5.263333
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0

# Row 40
$ws.Range("A40").Value = 1279
$ws.Range("B40").Value = "4:05 PM"
$ws.Range("C40").Value = 38365.04
$ws.Range("D40").Value = 38395.04
$ws.Range("E40").Value = 31.16
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0

# Row 41
$ws.Range("A41").Value = 1295
$ws.Range("B41").Value = "4:12 PM"
$ws.Range("C41").Value = 38838.2425
$ws.Range("D41").Value = 38868.2425
$ws.Range("E41").Value = 18.295
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0.25

# Row 42
$ws.Range("A42").Value = 1328
$ws.Range("B42").Value = "4:29 PM"
$ws.Range("C42").Value = 39833.97
$ws.Range("D42").Value = 39863.97
$ws.Range("E42").Value = 12.77
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0

# Row 43
$ws.Range("A43").Value = 1338
$ws.Range("B43").Value = "4:34 PM"
$ws.Range("C43").Value = 40124.17
$ws.Range("D43").Value = 40154.17
$ws.Range("E43").Value = 10.0375
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0.5

# Row 44
$ws.Range("A44").Value = 1353
$ws.Range("B44").Value = "4:42 PM"
$ws.Range("C44").Value = 40589.99
$ws.Range("D44").Value = 40619.99
$ws.Range("E44").Value = 4.08
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0

# Row 45
$ws.Range("A45").Value = 1361
$ws.Range("B45").Value = "4:45 PM"
$ws.Range("C45").Value = 40805.985
$ws.Range("D45").Value = 40835.985
$ws.Range("E45").Value = 17.225
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 3

# Row 46
$ws.Range("A46").Value = 1370
$ws.Range("B46").Value = "4:50 PM"
$ws.Range("C46").Value = 41074.77
$ws.Range("D46").Value = 41104.77
$ws.Range("E46").Value = 24.045
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 1

# Row 47
$ws.Range("A47").Value = 1390
$ws.Range("B47").Value = "5:00 PM"
$ws.Range("C47").Value = 41680.48
$ws.Range("D47").Value = 41710.48
$ws.Range("E47").Value = 38.76
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0.5

# Row 48
$ws.Range("A48").Value = 1409
$ws.Range("B48").Value = "5:09 PM"
$ws.Range("C48").Value = 42250.0775
$ws.Range("D48").Value = 42280.0775
$ws.Range("E48").Value = 4.4525
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0.25

# Row 49
$ws.Range("A49").Value = 1446
$ws.Range("B49").Value = "5:28 PM"
$ws.Range("C49").Value = 43374.21
$ws.Range("D49").Value = 43404.21
$ws.Range("E49").Value = 2.61
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0

# Row 50
$ws.Range("A50").Value = 1475
$ws.Range("B50").Value = "5:42 PM"
$ws.Range("C50").Value = 44227.06
$ws.Range("D50").Value = 44257.06
$ws.Range("E50").Value = 46.02
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 1.333333

# Row 51
$ws.Range("A51").Value = 1535
$ws.Range("B51").Value = "6:12 PM"
$ws.Range("C51").Value = 46024.31
$ws.Range("D51").Value = 46054.31
$ws.Range("E51").Value = 7.985
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0

# Give the two new header cells the same bold/border/centered style as the rest of row 1
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

